$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9-35 down to 10-36)
$ws.Rows("9:9").Insert()

# Populate the new row 9 with the latest weekly observation.
# Columns that stay identical to the former row 9 (now row 10) are copied
# across; only the date and the price columns change for this new record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = 'Vega Monumental Concepción'
$ws.Range("C9").Value = 'Bíobío'
$ws.Range("D9").Value = 44811
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112022
$ws.Range("G9").Value = 'Arveja Verde'
$ws.Range("H9").Value = 'Perfection'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 38000
$ws.Range("L9").Value = 40000
$ws.Range("M9").Value = 39000
$ws.Range("N9").Value = '$/malla 25 kilos'
$ws.Range("O9").Value = 'Provincia de Huasco'
$ws.Range("P9").Value = 1560
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = 'Hortaliza'

# Match the date-formatted style used by the other "Fecha" cells in column D
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
